$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 8
$cell = $ws.Cells.Item($newRow, 1)

$text = "meta features, baysian encoded features"
$address = "https://www.kaggle.com/c/bosch-production-line-performance/forums/t/24568/grid-search-for-xgb/136804"

# Put the new row's text in place first.
$cell.Value = $text

# Add the hyperlink pointing at the forum post (with its anchor/location),
# mirroring the other rows in the sheet.
$ws.Hyperlinks.Add($cell, $address, "post136804") | Out-Null

# Match the style used by the other hyperlink cells (A1:A7) exactly.
$cell.Style = $ws.Cells.Item(7, 1).Style

# Move/select the new last cell, like Excel would after entering data there.
$cell.Select() | Out-Null
